# Remove the small logo picture ("Picture 1") from every slide.
$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = $s.Shapes.Count; $j -ge 1; $j--) {
        $shape = $s.Shapes.Item($j)
        if ($shape.Name -eq "Picture 1") {
            $shape.Delete()
        }
    }
}
